$wb = $excel.ActiveWorkbook

# --- Worksheet references ---
$wsWeather = $wb.Worksheets.Item("Weather forecasts")
$wsPrices  = $wb.Worksheets.Item("Prices")

# --- Update values on "Weather forecasts" sheet ---
$wsWeather.Range("M2").Value = 0.613
$wsWeather.Range("N2").Value = 0.672
$wsWeather.Range("O2").Value = 0.565

# Move the selection on "Weather forecasts" sheet
$wsWeather.Range("Q9").Select()

# --- Update values on "Prices" sheet ---
$wsPrices.Range("G11:J11").Value = 75

# Make "Prices" the active sheet/tab and set its selection
$wsPrices.Activate()
$wsPrices.Range("F11:J11").Select()
